# Apply odds/score updates to Sheet1 as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2 updates ---
$ws.Range("G2").Value  = 1.57
$ws.Range("H2").Value  = 3.8
$ws.Range("I2").Value  = 6.5
$ws.Range("J2").Value  = 2.25
$ws.Range("K2").Value  = 2
$ws.Range("L2").Value  = 7.5
$ws.Range("M2").Value  = 1.11
$ws.Range("N2").Value  = 6.5
$ws.Range("W2").Value  = 4.75
$ws.Range("X2").Value  = 6
$ws.Range("AC2").Value = 6.5
$ws.Range("AD2").Value = 7.5
$ws.Range("AE2").Value = 26
$ws.Range("AF2").Value = 101
$ws.Range("AI2").Value = 21
$ws.Range("AJ2").Value = 81
$ws.Range("AL2").Value = 67
$ws.Range("AN2").Value = 3.25
$ws.Range("AO2").Value = 8.5
$ws.Range("AQ2").Value = 29
$ws.Range("AW2").Value = 8
$ws.Range("AZ2").Value = 201
$ws.Range("BA2").Value = 251

# --- Row 3 updates ---
$ws.Range("G3").Value  = 3.9
$ws.Range("H3").Value  = 2.8
$ws.Range("I3").Value  = 2.25
$ws.Range("J3").Value  = 4.75
$ws.Range("L3").Value  = 3.2
$ws.Range("Q3").Value  = 3.4
$ws.Range("R3").Value  = 1.33
$ws.Range("W3").Value  = 7
$ws.Range("X3").Value  = 17
$ws.Range("Y3").Value  = 15
$ws.Range("Z3").Value  = 41
$ws.Range("AB3").Value = 51
$ws.Range("AC3").Value = 5
$ws.Range("AF3").Value = 101
$ws.Range("AG3").Value = 5
$ws.Range("AH3").Value = 8.5
$ws.Range("AO3").Value = 26
$ws.Range("AW3").Value = 4

# --- Row 4/5 updates ---
$ws.Range("BD4").Value = 151
$ws.Range("BD5").Value = 126
